$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear out the old hyperlinks (the whole sheet collection; we'll re-add the one we keep) ---
$ws.Range("A1").Hyperlinks.Delete()

# --- Clear the old trailing "REST Protocol" filler rows (9-14 in col A only) ---
$ws.Range("A9:B14").ClearContents()

# --- Write cell values in the same order the original author introduced them,
#     so new shared-string entries land in the same append order. ---

# 1) A2 gets new (mangled, no-slash) local URL text
$ws.Range("A2").Value = "http://localhost:8080/MMarketWebServicegetAllPatronLoginHist"

# 2) B9 base production URL
$ws.Range("B9").Value = "http://apps-imeetem.com/MMarketWebService/"

# 3) A1 local server header
$ws.Range("A1").Value = "Man Market WebServices - LOCAL SERVER"

# 4) B1 production server header
$ws.Range("B1").Value = "Man Market WebServices - PRODUCTION SERVER"

# 5) B4 / B5 histId production URL (duplicated row)
$ws.Range("B4").Value = "http://apps-imeetem.com/MMarketWebService/getPatronLoginHistByHistId/{histId}"
$ws.Range("B5").Value = "http://apps-imeetem.com/MMarketWebService/getPatronLoginHistByHistId/{histId}"

# 6) B3 patronId production URL
$ws.Range("B3").Value = "http://apps-imeetem.com/MMarketWebService/getPatronLoginHistByPatronId/{patronId}"

# 7) B6 dateRange production URL
$ws.Range("B6").Value = "http://apps-imeetem.com/MMarketWebService/getPatronLoginHistByDateRange/{startDate}/{endDate}"

# 8) B7 updatePatronById production URL
$ws.Range("B7").Value = "http://apps-imeetem.com/MMarketWebService/updatePatronById/{patronId}/{lastLoginDate}"

# 9) B8 addPatronLoginHistRec production URL
$ws.Range("B8").Value = "http://apps-imeetem.com/MMarketWebService/addPatronLoginHistRec/{patronId}/{lastLoginDate}"

# 10) B2 getAllPatronLoginHist production URL (this is the real hyperlink cell)
$ws.Range("B2").Value = "http://apps-imeetem.com/MMarketWebService/getAllPatronLoginHist"

# --- C column: REST verbs ---
$ws.Range("C1").Value = "REST Protocol"
$ws.Range("C2").Value = "GET"
$ws.Range("C3").Value = "GET"
$ws.Range("C4").Value = "GET"
$ws.Range("C5").Value = "GET"
$ws.Range("C6").Value = "GET"
$ws.Range("C7").Value = "POST"
$ws.Range("C8").Value = "PUT"

# --- Header styling ---
$ws.Range("A1:B1").Style = "Normal"
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A1").Font.Size = 12
$ws.Range("B1").Font.Size = 12
$ws.Range("C1").Font.Size = 11

# --- Styling: plain (non-hyperlink) text cells A2:B9 except B2 ---
$ws.Range("A2:A8").Style = "Normal"
$ws.Range("B3:B9").Style = "Normal"

# --- Hyperlink: only B2 is a real hyperlink now ---
$ws.Hyperlinks.Add($ws.Range("B2"), "http://apps-imeetem.com/MMarketWebService/getAllPatronLoginHist")
$ws.Range("B2").Style = "Hyperlink"

# --- Column widths (closest achievable to 91.140625 / 97.85546875 / 27.42578125) ---
$ws.Columns.Item(1).ColumnWidth = 90.33
$ws.Columns.Item(2).ColumnWidth = 97
$ws.Columns.Item(3).ColumnWidth = 26.67

# --- Selection & view ---
$ws.Range("B14").Select()
